# "Generate Report for Handoff" - refresh the localization status report:
#  - flip the Status cells from "Handed back: in sync with en-US" to "Ready for handoff"
#  - bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#  - shrink the now-narrower Status columns to fit the shorter text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: zh-cn (E2) and de-de (F2) status columns, plus generate date (G2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-09-06 05:16:56"

# zh-cn detail sheet: status (C2) and latest handoff datetime (H2)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-09-06 05:16:43"

# de-de detail sheet: status (C2) and latest handoff datetime (H2)
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-09-06 05:16:56"

# Re-fit the Status columns now that the text is shorter
$overview.Columns("E:F").ColumnWidth = 16.33
$zhcn.Columns("C:C").ColumnWidth = 16.33
$dede.Columns("C:C").ColumnWidth = 16.33
